# "añadidos datos de 125 y 250 a 2002.xlsx y sus inputs"
# The main_list sheet ("Hoja1") previously listed every year's data file
# (2002-2023). This edit trims it back down to just the header row plus the
# 2002 entry (rows 3-23, i.e. years 2003-2023, are removed).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate() | Out-Null

# Remove rows 3 through 23, leaving only the header (row 1) and the 2002
# data row (row 2).
$ws.Range("A3:B23").EntireRow.Delete() | Out-Null

# Match the post-edit selection recorded for this sheet.
$ws.Range("A3:B23").Select() | Out-Null
